$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the "MCT-1A-Metrologia" class block from Monday (column B) to Thursday (column E)
# for the 7:50, 8:40, 9:50 and 10:40 slots, and add it to the 11:30 Thursday slot.

$ws.Range("B3").Value = "-"
$ws.Range("B4").Value = "-"
$ws.Range("E4").Value = "[-, -, -, 'MCT-1A-Metrologia']"
$ws.Range("B6").Value = "-"
$ws.Range("E6").Value = "[-, -, -, 'MCT-1A-Metrologia']"
$ws.Range("B7").Value = "-"
$ws.Range("E7").Value = "[-, -, -, 'MCT-1A-Metrologia']"
$ws.Range("E8").Value = "[-, -, -, 'MCT-1A-Metrologia']"
